# Scheduled-runner price/profit refresh across the Sheets workbook.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# for a handful of leve rows on each crafting-class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5197.3184
$ws.Range("I40").Value = 5993.769
$ws.Range("J40").Value = 4046.889
$ws.Range("K40").Value = 5993.769
$ws.Range("L40").Value = 4046.889
$ws.Range("M40").Value = -5818.769
$ws.Range("N40").Value = -4396.889

$ws.Range("H88").Value = 5127.7407
$ws.Range("J88").Value = 4950.263
$ws.Range("L88").Value = 4950.263
$ws.Range("N88").Value = -5762.263

$ws.Range("H91").Value = 5127.7407
$ws.Range("J91").Value = 4950.263
$ws.Range("L91").Value = 4950.263
$ws.Range("N91").Value = -7758.263

$ws.Range("H125").Value = 1610.875
$ws.Range("I125").Value = 1099
$ws.Range("J125").Value = 1781.5
$ws.Range("K125").Value = 9891
$ws.Range("L125").Value = 16033.5
$ws.Range("M125").Value = -7431
$ws.Range("N125").Value = -20953.5

$ws.Range("H126").Value = 60000
$ws.Range("J126").Value = 60000
$ws.Range("L126").Value = 60000
$ws.Range("N126").Value = -69880

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4901.049
$ws.Range("I32").Value = 5179.4473
$ws.Range("K32").Value = 5179.4473
$ws.Range("M32").Value = -4892.4473

$ws.Range("H132").Value = 1497.0769
$ws.Range("I132").Value = 1544.75
$ws.Range("K132").Value = 4634.25
$ws.Range("M132").Value = -2104.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5763.24
$ws.Range("I31").Value = 10773.2
$ws.Range("J31").Value = 2423.2666
$ws.Range("K31").Value = 10773.2
$ws.Range("L31").Value = 2423.2666
$ws.Range("M31").Value = -10478.2
$ws.Range("N31").Value = -3013.2666

$ws.Range("H34").Value = 5763.24
$ws.Range("I34").Value = 10773.2
$ws.Range("J34").Value = 2423.2666
$ws.Range("K34").Value = 10773.2
$ws.Range("L34").Value = 2423.2666
$ws.Range("M34").Value = -10571.2
$ws.Range("N34").Value = -2827.2666

$ws.Range("H118").Value = 49998.75
$ws.Range("J118").Value = 49998.75
$ws.Range("L118").Value = 49998.75
$ws.Range("N118").Value = -53312.75

$ws.Range("H122").Value = 304776.66
$ws.Range("J122").Value = 4676.2
$ws.Range("L122").Value = 14028.6
$ws.Range("N122").Value = -18928.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 549.3333
$ws.Range("I8").Value = 549.3333
$ws.Range("K8").Value = 1647.9999
$ws.Range("M8").Value = -1508.9999

$ws.Range("H23").Value = 232
$ws.Range("J23").Value = 232
$ws.Range("L23").Value = 696
$ws.Range("N23").Value = -1166

$ws.Range("H37").Value = 405885.4
$ws.Range("J37").Value = 405885.4
$ws.Range("L37").Value = 1217656.2
$ws.Range("N37").Value = -1217880.2

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H113").Value = 766.8570999999999
$ws.Range("I113").Value = 971.1667
$ws.Range("J113").Value = 613.625
$ws.Range("K113").Value = 2913.5001
$ws.Range("L113").Value = 1840.875
$ws.Range("M113").Value = -743.5001000000002
$ws.Range("N113").Value = -6180.875

$ws.Range("H121").Value = 2881.75
$ws.Range("J121").Value = 5121.4165
$ws.Range("L121").Value = 15364.2495
$ws.Range("N121").Value = -17984.2495

$ws.Range("H129").Value = 1591.1428
$ws.Range("I129").Value = 1075.1428
$ws.Range("J129").Value = 2107.1428
$ws.Range("K129").Value = 3225.4284
$ws.Range("L129").Value = 6321.428400000001
$ws.Range("M129").Value = 1774.5716
$ws.Range("N129").Value = -16321.4284

$ws.Range("H131").Value = 1669090.4
$ws.Range("I131").Value = 2000788.5
$ws.Range("J131").Value = 10600
$ws.Range("K131").Value = 6002365.5
$ws.Range("L131").Value = 31800
$ws.Range("M131").Value = -5997325.5
$ws.Range("N131").Value = -41880

$ws.Range("H137").Value = 1750
$ws.Range("I137").Value = 1750
$ws.Range("K137").Value = 5250
$ws.Range("M137").Value = -150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 188843.33
$ws.Range("I80").Value = 560555
$ws.Range("K80").Value = 560555
$ws.Range("M80").Value = -559557

$ws.Range("H83").Value = 188843.33
$ws.Range("I83").Value = 560555
$ws.Range("K83").Value = 2802775
$ws.Range("M83").Value = -2797783

$ws.Range("H97").Value = 12740.357
$ws.Range("I97").Value = 4555.75
$ws.Range("K97").Value = 4555.75
$ws.Range("M97").Value = -4059.75

$ws.Range("H124").Value = 57400
$ws.Range("J124").Value = 57400
$ws.Range("L124").Value = 57400
$ws.Range("N124").Value = -67220

$ws.Range("H126").Value = 3441.3635
$ws.Range("J126").Value = 4353
$ws.Range("L126").Value = 13059
$ws.Range("N126").Value = -17999

$ws.Range("H135").Value = 83696.75999999999
$ws.Range("J135").Value = 83696.75999999999
$ws.Range("L135").Value = 83696.75999999999
$ws.Range("N135").Value = -93836.75999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2525.4783
$ws.Range("J22").Value = 2000.25
$ws.Range("L22").Value = 2000.25
$ws.Range("N22").Value = -2590.25

$ws.Range("H27").Value = 2525.4783
$ws.Range("J27").Value = 2000.25
$ws.Range("L27").Value = 2000.25
$ws.Range("N27").Value = -2214.25

$ws.Range("H61").Value = 5666
$ws.Range("I61").Value = 5666
$ws.Range("K61").Value = 5666
$ws.Range("M61").Value = -5464

$ws.Range("H113").Value = 5666
$ws.Range("I113").Value = 5666
$ws.Range("K113").Value = 5666
$ws.Range("M113").Value = -3496

$ws.Range("H122").Value = 4699.2
$ws.Range("I122").Value = 2864.8333
$ws.Range("K122").Value = 8594.499899999999
$ws.Range("M122").Value = -6144.499899999999

$ws.Range("H132").Value = 5190.2
$ws.Range("I132").Value = 5378
$ws.Range("K132").Value = 16134
$ws.Range("M132").Value = -13604

$ws.Range("H137").Value = 80000
$ws.Range("J137").Value = 80000
$ws.Range("L137").Value = 80000
$ws.Range("N137").Value = -90200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1338.8235
$ws.Range("J122").Value = 1432.3334
$ws.Range("L122").Value = 4297.0002
$ws.Range("N122").Value = -9197.0002

$ws.Range("H126").Value = 3130.3333
$ws.Range("I126").Value = 2932.1177
$ws.Range("K126").Value = 8796.3531
$ws.Range("M126").Value = -6326.3531
